$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the three columns that are no longer part of the report:
#    Téléphone, Modalité de réception, Cause de retour.
$ws.Range("H1:J1").EntireColumn.Delete()

# 2. Insert a new column in front of the old "État" column (now shifted to L)
#    to host the new "Valeur" sub-header.
$ws.Range("L1").EntireColumn.Insert()

# 3. Merge the new header cell and copy formatting from a sibling sub-header
#    cell so it matches the surrounding "Informations Articles" styling.
$ws.Range("L2:L3").Merge()
$ws.Range("H2").Copy()
$ws.Range("L2:L3").PasteSpecial(-4122)
$ws.Range("Z1").Clear()

# 4. Update the header labels that changed wording.
$ws.Range("A1").Value = "Bon Livraison"
$ws.Range("B1").Value = "Wilaya"
$ws.Range("G1").Value = "Probleme"
$ws.Range("L2").Value = "Valeur"

# Helper: write a value as TEXT (avoids Excel auto-converting numeric-looking
# or date-looking strings) while preserving the destination cell's existing
# style, by staging the text in a scratch cell formatted as Text and pasting
# only the value back onto the target cell.
function Set-TextValue($addr, $text) {
    $ws.Range("Z1").NumberFormat = "@"
    $ws.Range("Z1").Value = $text
    $ws.Range("Z1").Copy()
    $ws.Range($addr).PasteSpecial(-4163)
    $ws.Range("Z1").Clear()
}

# 5. Refresh the data row with the new sample values.
Set-TextValue "A4" "7257845278"
$ws.Range("B4").Value = "Alger"
Set-TextValue "C4" "12/08/2024"
Set-TextValue "D4" "01/08/2024"
$ws.Range("E4").Value = "LAIB  HAMID"
Set-TextValue "F4" "0635515554"
$ws.Range("G4").Value = "erreurLivraison"
$ws.Range("H4").Value = "OPT001CH"
$ws.Range("I4").Value = "INTERRUPTEUR SIMPLE ALLUMAGE"
Set-TextValue "J4" "25"
Set-TextValue "K4" ""
$ws.Range("L4").Value = "11000.DA"
$ws.Range("M4").Value = "En plus"

# 6. Add the trailing padding column (N) that keeps the grid 14 columns wide,
#    matching the plain centred style Excel produced for spacer cells here.
$ws.Range("N1").HorizontalAlignment = -4108
$ws.Range("N1").VerticalAlignment = -4108
$ws.Range("N4").HorizontalAlignment = -4108
$ws.Range("N4").VerticalAlignment = -4108
